$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.654
$ws.Range("C5").Value = -12.156
$ws.Range("E7").Value = 13.078
$ws.Range("C9").Value = -11.671
$ws.Range("C11").Value = -12.628
$ws.Range("E11").Value = 12.942
$ws.Range("A21").Value = -21.261
$ws.Range("C21").Value = -12.442
$ws.Range("E21").Value = 13.202
$ws.Range("A23").Value = -21.654
$ws.Range("A25").Value = -21.964

$wb.Save()
